$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Values (written in the order that reproduces the original shared-string layout) ---

# Header row + first two data rows, columns A-D
$ws.Range("A17").Value = "UC"
$ws.Range("B17").Value = "Documentation"
$ws.Range("C17").Value = "Coding"
$ws.Range("D17").Value = "Testing"

$ws.Range("A18").Value = "Detect Face"
$ws.Range("B18").Value = "2h"
$ws.Range("C18").Value = "18h"
$ws.Range("D18").Value = "3h"

$ws.Range("A19").Value = "Label Image"
$ws.Range("B19").Value = "1h"
$ws.Range("C19").Value = "25h"
$ws.Range("D19").Value = "5h"

# Header row + first two data rows, columns E-F (Total / FP)
$ws.Range("E17").Value = "Total"
$ws.Range("F17").Value = "FP"

$ws.Range("E18").Value = "23h"
$ws.Range("F18").Formula = "=B3"

$ws.Range("E19").Value = "31h"
$ws.Range("F19").Formula = "=B2"

# Remaining three data rows, all columns A-F
$ws.Range("A20").Value = "Upload Image"
$ws.Range("B20").Value = "1h"
$ws.Range("C20").Value = "5h"
$ws.Range("D20").Value = "0,5h"
$ws.Range("E20").Value = "6,5h"
$ws.Range("F20").Formula = "=B6"

$ws.Range("A21").Value = "Delete Image"
$ws.Range("B21").Value = "1h"
$ws.Range("C21").Value = "6h"
$ws.Range("D21").Value = "1h"
$ws.Range("E21").Value = "8h"
$ws.Range("F21").Formula = "=B5"

$ws.Range("A22").Value = "Browse Image"
$ws.Range("B22").Value = "1h"
$ws.Range("C22").Value = "12h"
$ws.Range("D22").Value = "1h"
$ws.Range("E22").Value = "14h"
$ws.Range("F22").Formula = "=B4"

# --- Formatting ---

# Row 17: bold 14pt header style (copy existing bold-header look), then recolor the fill
$ws.Range("A1:C1").Copy()
$ws.Range("A17:F17").PasteSpecial(-4122)
$ws.Range("A17:F17").Interior.Color = 0xA6A6A6

# Rows 18/20/22 use the plain banded style, rows 19/21 use the shaded banded style
$ws.Range("A2:C2").Copy()
$ws.Range("A18:F18").PasteSpecial(-4122)

$ws.Range("A3:C3").Copy()
$ws.Range("A19:F19").PasteSpecial(-4122)

$ws.Range("A2:C2").Copy()
$ws.Range("A20:F20").PasteSpecial(-4122)

$ws.Range("A3:C3").Copy()
$ws.Range("A21:F21").PasteSpecial(-4122)

$ws.Range("A2:C2").Copy()
$ws.Range("A22:F22").PasteSpecial(-4122)

# Restore the original selection location as seen in the authored workbook
$ws.Range("C25").Select() | Out-Null
